$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.1360299999124
$ws.Range("C2").Value = 381
$ws.Range("E2").Value = [double]"-7.72051916934829E-07"
$ws.Range("F2").Value = 0.2417821421664151
$ws.Range("G2").Value = 3695.204786682511
$ws.Range("H2").Value = 0.5719853491229088
$ws.Range("B3").Value = 21.84980946997159
$ws.Range("C3").Value = 5
$ws.Range("F3").Value = 0.2821717707289705
$ws.Range("G3").Value = 3773.474156720286
$ws.Range("H3").Value = 0.5790369448021434
$ws.Range("B4").Value = 22.55953039997485
$ws.Range("E4").Value = [double]"3.86022916934829E-07"
$ws.Range("F4").Value = 0.3026597777944491
$ws.Range("G4").Value = 3914.675038632136
$ws.Range("H4").Value = 0.576281049572319
$ws.Range("B5").Value = 23.22413391990703
$ws.Range("F5").Value = 0.3158736445565552
$ws.Range("G5").Value = 4076.780433409805
$ws.Range("H5").Value = 0.5696684994262113
$ws.Range("B6").Value = 23.9552675399555
$ws.Range("F6").Value = 0.3200592976877066
$ws.Range("G6").Value = 4257.647815580628
$ws.Range("H6").Value = 0.5626408894670084
$ws.Range("E7").Value = 0.5770523579999891
$ws.Range("F7").Value = 0.3655146931944009
$ws.Range("G7").Value = 4455.774820691072
$ws.Range("H7").Value = 0.5558253842847056
$ws.Range("B8").Value = 25.65698091995749
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 14.662803998
$ws.Range("F8").Value = 0.3661127750801957
$ws.Range("G8").Value = 4597.276055748865
$ws.Range("H8").Value = 0.5580909349107631
$ws.Range("B9").Value = 26.73718129996622
$ws.Range("E9").Value = 32.20124931477594
$ws.Range("F9").Value = 0.386180324527109
$ws.Range("G9").Value = 4675.811406468223
$ws.Range("H9").Value = 0.5718190700116709
$ws.Range("B10").Value = 27.09405198996807
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 86.41899574346111
$ws.Range("F10").Value = 0.3865908209481895
$ws.Range("G10").Value = 4757.105457017769
$ws.Range("H10").Value = 0.5695491141571905
$ws.Range("B11").Value = 27.45095816996313
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 17
$ws.Range("E11").Value = 116.3109139457456
$ws.Range("F11").Value = 0.3780578445903078
$ws.Range("G11").Value = 4875.291586712862
$ws.Range("H11").Value = 0.5630628995561634
$ws.Range("B12").Value = 27.81457391999498
$ws.Range("C12").Value = 21
$ws.Range("D12").Value = 34
$ws.Range("E12").Value = 274.1592631083323
$ws.Range("F12").Value = 0.3353575826160597
$ws.Range("G12").Value = 5065.234685496795
$ws.Range("H12").Value = 0.5491270522891271
$ws.Range("B13").Value = 28.16623758996121
$ws.Range("C13").Value = 2.99999945548237
$ws.Range("D13").Value = 33
$ws.Range("E13").Value = 260.7024028676158
$ws.Range("F13").Value = 0.3236666966834131
$ws.Range("G13").Value = 5231.405173578296
$ws.Range("H13").Value = 0.5384067311822346
$ws.Range("B14").Value = 28.44839399996198
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 30
$ws.Range("E14").Value = 254.6208081453192
$ws.Range("F14").Value = 0.3237615998969312
$ws.Range("G14").Value = 5339.21466325712
$ws.Range("H14").Value = 0.5328198207825455
$ws.Range("B15").Value = 28.49488944995075
$ws.Range("C15").Value = 5
$ws.Range("E15").Value = 292.2713247900767
$ws.Range("F15").Value = 0.3110630587953992
$ws.Range("G15").Value = 5417.962626924868
$ws.Range("H15").Value = 0.5259336656983162
$ws.Range("B16").Value = 28.52573986994972
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 31
$ws.Range("E16").Value = 210.8172857962811
$ws.Range("F16").Value = 0.2858120738669336
$ws.Range("G16").Value = 5514.972687740342
$ws.Range("H16").Value = 0.5172417251197234
$ws.Range("B17").Value = 28.5455398399759
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 24
$ws.Range("E17").Value = 176.3400327650926
$ws.Range("F17").Value = 0.2755856237523093
$ws.Range("G17").Value = 5566.585602776334
$ws.Range("H17").Value = 0.5128015964712519
$ws.Range("B18").Value = 28.52921014994949
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 21
$ws.Range("E18").Value = 163.4850990038484
$ws.Range("F18").Value = 0.2754038949529412
$ws.Range("G18").Value = 5685.15439485233
$ws.Range("H18").Value = 0.5018194435630718
$ws.Range("B19").Value = 28.49762812995008
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = 145.9223132966955
$ws.Range("F19").Value = 0.2776940333582884
$ws.Range("G19").Value = 5738.589918072874
$ws.Range("H19").Value = 0.4965963509642125
$ws.Range("B20").Value = 28.27933570994956
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = 71.38904654452021
$ws.Range("F20").Value = 0.2769542029011573
$ws.Range("G20").Value = 5724.266261662869
$ws.Range("H20").Value = 0.4940255120441333
$ws.Range("B21").Value = 28.04663652994937
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 7
$ws.Range("E21").Value = 39.33091134198768
$ws.Range("F21").Value = 0.2579195429379323
$ws.Range("G21").Value = 5703.459761074591
$ws.Range("H21").Value = 0.4917477759966714
$ws.Range("B22").Value = 27.78277931994786
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = 30.12405280171151
$ws.Range("F22").Value = 0.247987002089462
$ws.Range("G22").Value = 5708.644878901238
$ws.Range("H22").Value = 0.486679061481493
$ws.Range("B23").Value = 27.493555829949
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 10.67804252524875
$ws.Range("F23").Value = 0.2475432192726494
$ws.Range("G23").Value = 5653.400613769043
$ws.Range("H23").Value = 0.4863189026970341
$ws.Range("B24").Value = 27.18263807994779
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 0.153541385
$ws.Range("F24").Value = 0.2475445761891555
$ws.Range("G24").Value = 5595.526865739282
$ws.Range("H24").Value = 0.4857922896659422
$ws.Range("B25").Value = 26.78460247004762
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 48.588705508
$ws.Range("F25").Value = 0.2439671194956984
$ws.Range("G25").Value = 5503.665282162981
$ws.Range("H25").Value = 0.4866684490580261
$ws.Range("B26").Value = 26.37472910995696
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 127.5728948626793
$ws.Range("F26").Value = 0.2406107521044039
$ws.Range("G26").Value = 5392.651014743664
$ws.Range("H26").Value = 0.489086518631517
$ws.Range("B27").Value = 25.98391770995499
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 19
$ws.Range("E27").Value = 204.8591429804337
$ws.Range("F27").Value = 0.235379673115279
$ws.Range("G27").Value = 5347.939877339449
$ws.Range("H27").Value = 0.4858677978048204
$ws.Range("B28").Value = 25.60600485995852
$ws.Range("C28").Value = 11.00000009033469
$ws.Range("E28").Value = 266.1911040768323
$ws.Range("F28").Value = 0.2246178450275679
$ws.Range("G28").Value = 5275.456998061773
$ws.Range("H28").Value = 0.4853798423409817
$ws.Range("B29").Value = 25.22753590995961
$ws.Range("C29").Value = 14
$ws.Range("D29").Value = 36
$ws.Range("E29").Value = 266.7352401030083
$ws.Range("F29").Value = 0.2252530815462508
$ws.Range("G29").Value = 5294.818652612466
$ws.Range("H29").Value = 0.4764570340385941
